# Add three new price columns ("Precio Residencial", "Precio Comercial",
# "Precio Distribuidor") right after the existing "Precio" column, and
# refresh the sheet's title text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Precio" lives in column E; the new fields are inserted immediately
# after it, pushing the remaining headers (Nombre del material, Color,
# Ancho, Composicion, Flamabilidad, Minimos/Multiplos de venta, Tamano
# de rollo) three columns to the right.
$ws.Columns("F:H").Insert()

$ws.Range("F2").Value = "Precio Residencial"
$ws.Range("G2").Value = "Precio Comercial"
$ws.Range("H2").Value = "Precio Distribuidor"

# The merged title banner auto-extends with the inserted columns; just
# refresh its text to the updated example-file name (no more "(1)" suffix).
$ws.Range("A1").Value = "ejemplo_formato_carga_masiva_productos_tapices_-_formato_carga_masiva_productos_tapices"
